# Refresh the cryptos price/volume table (columns D and E, rows 2-51).
# Source cells are plain text (e.g. "61.845.77", "0.999", "  -0.83%  ") rather
# than real numbers, so for values that Excel would otherwise auto-parse as a
# number (e.g. "1.00" -> 1) we write them with a leading apostrophe to force
# text entry, then reset the cell style to "Normal" so the quote-prefix flag
# does not leave a stray number-format on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.924.86"
$ws.Range("D3").Value = "3.413.77"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'410.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "'129.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").Value = "'0.639"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.20%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.738"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.41%  "
$ws.Range("D10").Value = "'0.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").Value = "'43.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").Value = "'0.0000225"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +37.22%  "
$ws.Range("D13").Value = "'9.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.08%  "
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "'21.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.08%  "
$ws.Range("D16").Value = "3.951.30"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "3.404.62"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "'12.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.50%  "
$ws.Range("E19").Value = "  +6.30%  "
$ws.Range("D20").Value = "61.905.30"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").Value = "'467.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +47.65%  "
$ws.Range("D22").Value = "'92.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.02%  "
$ws.Range("D23").Value = "'3.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "'13.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").Value = "'3.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.27%  "
$ws.Range("D26").Value = "'33.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.31%  "
$ws.Range("D27").Value = "'9.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.73%  "
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").Value = "'7.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("E31").Value = "  +4.53%  "
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").Value = "'42.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.26%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'0.0506"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.49%  "
$ws.Range("D37").Value = "'53.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.87%  "
$ws.Range("D38").Value = "'0.997"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "'0.137"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.03%  "
$ws.Range("D40").Value = "'3.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").Value = "'0.321"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "'4.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.93%  "
$ws.Range("D44").Value = "'144.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +16.53%  "
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").Value = "'16.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "'0.151"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +20.87%  "
$ws.Range("D49").Value = "'22.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.11%  "
$ws.Range("D50").Value = "'2.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.60%  "
$ws.Range("D51").Value = "3.753.91"
$ws.Range("E51").Value = "  -0.76%  "
